$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "wj" column (G) and the "wij" column (which becomes H after the
# first deletion shifts everything left one place).
$ws.Columns("G").Delete()
$ws.Columns("H").Delete()

# The remaining "omega_ij" column (now I) holds recalculated values that are
# no longer simply a shifted copy of the old K column, so set them directly.
$ws.Range("I2").Value = 0.3029128375207381
$ws.Range("I3").Value = 0.2645643692038547
$ws.Range("I4").Value = 0.5829401856184298
